# "added final result calcs after adjusting ML labels"
#
# Sheet2 holds the "DJI_pos_50" table (date / Close / Label) that is
# normally filled by the Power Query connection, with Label (column C)
# acting as the ML-predicted buy/hold flag used by the D/E/H/I/J trading
# calcs below it. Here the ML labels for a handful of dates were revised
# by hand, overwriting the query's cached Label values; everything else
# (D, E, H2, H3, I2/I3, J3 ...) is formula-driven off column C and just
# recalculates from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$labelChanges = @{
    2  = 1
    3  = 0
    5  = 1
    6  = 1
    10 = 0
    13 = 0
    16 = 1
    17 = 0
    18 = 0
    23 = 0
    25 = 0
    26 = 1
    28 = 0
    31 = 1
    32 = 1
    34 = 1
    38 = 0
    42 = 1
    43 = 0
    46 = 1
    48 = 1
    49 = 0
    50 = 1
    51 = 1
}

foreach ($row in $labelChanges.Keys) {
    $ws.Range("C$row").Value = $labelChanges[$row]
}
